# Update the "About" sheet: add a region label (Oregon) next to the
# existing title cell, and bump the "last updated" date stamp.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("BIEfIE")

$ws1.Range("B1").Value = "Oregon"
$ws1.Range("C1").Value = 44840

# Make the BIEfIE (control-lever) sheet the active tab, with B3 selected,
# matching the saved workbook/view state in the updated file.
$ws2.Activate()
$ws2.Range("B3").Select()
